# Scheduled data-runner refresh: updates market-price-derived columns
# (currentAveragePrice*, LevePrice*, LeveProfit*) on several sheets with
# freshly pulled numbers. Plain value writes, no formulas involved.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 10170.5
$ws.Range("J51").Value = 8000
$ws.Range("L51").Value = 8000
$ws.Range("N51").Value = -8968
$ws.Range("H58").Value = 5885.9375
$ws.Range("I58").Value = 423.66666
$ws.Range("J58").Value = 9163.299999999999
$ws.Range("K58").Value = 1270.99998
$ws.Range("L58").Value = 27489.9
$ws.Range("M58").Value = -1120.99998
$ws.Range("N58").Value = -27789.9
$ws.Range("H64").Value = 6287.1914
$ws.Range("I64").Value = 4532.6665
$ws.Range("J64").Value = 6406.8184
$ws.Range("K64").Value = 4532.6665
$ws.Range("L64").Value = 6406.8184
$ws.Range("M64").Value = -4284.6665
$ws.Range("N64").Value = -6902.8184
$ws.Range("H67").Value = 6287.1914
$ws.Range("I67").Value = 4532.6665
$ws.Range("J67").Value = 6406.8184
$ws.Range("K67").Value = 4532.6665
$ws.Range("L67").Value = 6406.8184
$ws.Range("M67").Value = -3674.6665
$ws.Range("N67").Value = -8122.8184
$ws.Range("H74").Value = 13298.385
$ws.Range("J74").Value = 21665
$ws.Range("L74").Value = 21665
$ws.Range("N74").Value = -23537
$ws.Range("H77").Value = 13298.385
$ws.Range("J77").Value = 21665
$ws.Range("L77").Value = 108325
$ws.Range("N77").Value = -117685

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11975.85
$ws.Range("I32").Value = 12185.368
$ws.Range("K32").Value = 12185.368
$ws.Range("M32").Value = -11898.368
$ws.Range("H61").Value = 4649.643
$ws.Range("I61").Value = 4315
$ws.Range("K61").Value = 4315
$ws.Range("M61").Value = -4103
$ws.Range("H74").Value = 2071.724
$ws.Range("I74").Value = 2327.2104
$ws.Range("J74").Value = 1586.3
$ws.Range("K74").Value = 2327.2104
$ws.Range("L74").Value = 1586.3
$ws.Range("M74").Value = -1453.2104
$ws.Range("N74").Value = -3334.3
$ws.Range("H77").Value = 2071.724
$ws.Range("I77").Value = 2327.2104
$ws.Range("J77").Value = 1586.3
$ws.Range("K77").Value = 11636.052
$ws.Range("L77").Value = 7931.5
$ws.Range("M77").Value = -7268.052
$ws.Range("N77").Value = -16667.5
$ws.Range("H122").Value = 3387.2632
$ws.Range("J122").Value = 5018.8945
$ws.Range("L122").Value = 15056.6835
$ws.Range("N122").Value = -19956.6835
$ws.Range("H136").Value = 4649.643
$ws.Range("I136").Value = 4315
$ws.Range("K136").Value = 12945
$ws.Range("M136").Value = -10395

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 30572.316
$ws.Range("I134").Value = 5012.52
$ws.Range("K134").Value = 15037.56
$ws.Range("M134").Value = -12502.56

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 6306.643
$ws.Range("I99").Value = 5471.2856
$ws.Range("J99").Value = 7142
$ws.Range("K99").Value = 5471.2856
$ws.Range("L99").Value = 7142
$ws.Range("M99").Value = -3973.2856
$ws.Range("N99").Value = -10138
$ws.Range("H126").Value = 6306.643
$ws.Range("I126").Value = 5471.2856
$ws.Range("J126").Value = 7142
$ws.Range("K126").Value = 16413.8568
$ws.Range("L126").Value = 21426
$ws.Range("M126").Value = -13943.8568
$ws.Range("N126").Value = -26366
$ws.Range("H134").Value = 348170.1
$ws.Range("I134").Value = 3231.5386
$ws.Range("K134").Value = 9694.6158
$ws.Range("M134").Value = -7159.6158

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H74").Value = 4358
$ws.Range("J74").Value = 6315
$ws.Range("L74").Value = 18945
$ws.Range("N74").Value = -21067
$ws.Range("H77").Value = 4358
$ws.Range("J77").Value = 6315
$ws.Range("L77").Value = 56835
$ws.Range("N77").Value = -67443
$ws.Range("H107").Value = 59396.277
$ws.Range("I107").Value = 1116.8572
$ws.Range("J107").Value = 96483.17999999999
$ws.Range("K107").Value = 3350.5716
$ws.Range("L107").Value = 289449.54
$ws.Range("M107").Value = -1430.5716
$ws.Range("N107").Value = -293289.54
$ws.Range("H113").Value = 1059809.4
$ws.Range("J113").Value = 1917.7273
$ws.Range("L113").Value = 5753.1819
$ws.Range("N113").Value = -10093.1819
$ws.Range("H122").Value = 2053.4
$ws.Range("I122").Value = 1088.4
$ws.Range("J122").Value = 2535.9
$ws.Range("K122").Value = 9795.6
$ws.Range("L122").Value = 22823.1
$ws.Range("M122").Value = -7345.6
$ws.Range("N122").Value = -27723.1
$ws.Range("H126").Value = 4262.25
$ws.Range("I126").Value = 4262.25
$ws.Range("K126").Value = 12786.75
$ws.Range("M126").Value = -7846.75
$ws.Range("H129").Value = 63649.25
$ws.Range("J129").Value = 168961
$ws.Range("L129").Value = 506883
$ws.Range("N129").Value = -516883
$ws.Range("H137").Value = 2863.138
$ws.Range("I137").Value = 2544.4
$ws.Range("J137").Value = 3571.4443
$ws.Range("K137").Value = 7633.200000000001
$ws.Range("L137").Value = 10714.3329
$ws.Range("M137").Value = -2533.200000000001
$ws.Range("N137").Value = -20914.3329
$ws.Range("H139").Value = 4796.5405
$ws.Range("I139").Value = 1074.8572
$ws.Range("J139").Value = 9681.25
$ws.Range("K139").Value = 3224.5716
$ws.Range("L139").Value = 29043.75
$ws.Range("M139").Value = 1915.4284
$ws.Range("N139").Value = -39323.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H136").Value = 27193.932
$ws.Range("J136").Value = 27193.932
$ws.Range("L136").Value = 81581.796
$ws.Range("N136").Value = -86681.796

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 113102.78
$ws.Range("I68").Value = 2275
$ws.Range("J68").Value = 501000
$ws.Range("K68").Value = 2275
$ws.Range("L68").Value = 501000
$ws.Range("M68").Value = -1526
$ws.Range("N68").Value = -502498
$ws.Range("H71").Value = 113102.78
$ws.Range("I71").Value = 2275
$ws.Range("J71").Value = 501000
$ws.Range("K71").Value = 11375
$ws.Range("L71").Value = 2505000
$ws.Range("M71").Value = -7631
$ws.Range("N71").Value = -2512488
$ws.Range("H132").Value = 3702.4827
$ws.Range("I132").Value = 3606.3704
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 10819.1112
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -8289.111199999999
$ws.Range("N132").Value = -20060
$ws.Range("H136").Value = 257152.33
$ws.Range("I136").Value = 481984.34
$ws.Range("K136").Value = 1445953.02
$ws.Range("M136").Value = -1443403.02

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 3666.6667
$ws.Range("I14").Value = 3666.6667
$ws.Range("K14").Value = 3666.6667
$ws.Range("M14").Value = -3498.6667
$ws.Range("H81").Value = 5614.5864
$ws.Range("I81").Value = 2373.2222
$ws.Range("K81").Value = 4746.4444
$ws.Range("M81").Value = -3685.4444
$ws.Range("H84").Value = 5614.5864
$ws.Range("I84").Value = 2373.2222
$ws.Range("K84").Value = 23732.222
$ws.Range("M84").Value = -18428.222
$ws.Range("H132").Value = 74145.664
$ws.Range("I132").Value = 7743.636
$ws.Range("K132").Value = 23230.908
$ws.Range("M132").Value = -20700.908
$ws.Range("H136").Value = 89247
$ws.Range("I136").Value = 22897.762
$ws.Range("K136").Value = 68693.28599999999
$ws.Range("M136").Value = -66143.28599999999
